$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that change per row (row 2..17), in order
$cols = @("E","G","H","I","J","K","M","N","O","P","Q","R","S","T")

# New values for each data row (row indices 2..17), one array per row in $cols order
$rowData = @(
    @(3, 2.861854666666666, 8.585564, 0.05512076040054963, 0.05512076040054963, 3, 7.248785666666667, 21.746357, 0.07891374419744837, 0.07891374419744837, 20.74497108781645, 186.704739790348, 0.004349785586217816, 0.004349785586217816),
    @(3, 2.861854666666666, 8.585564, 0.05512076040054963, 0.05512076040054963, 3, 37.82684066666667, 113.480522, 0.4118010609547572, 0.4118010609547572, 108.2549204871564, 974.294284384408, 0.02269878761357931, 0.02269878761357931),
    @(3, 2.861854666666666, 8.585564, 0.05512076040054963, 0.05512076040054963, 3, 38.20927633333334, 114.627829, 0.4159644383477588, 0.4159644383477588, 109.3493957845062, 984.144562060556, 0.02292827614131601, 0.02292827614131601),
    @(3, 2.861854666666666, 8.585564, 0.05512076040054963, 0.05512076040054963, 3, 8.572171666666666, 25.716515, 0.09332075650003555, 0.09332075650003555, 24.53230948771778, 220.79078538946, 0.005143911059436494, 0.005143911059436494),
    @(3, 19.097892, 57.293676, 0.3678350062107417, 0.3678350062107418, 3, 7.248785666666667, 21.746357, 0.07891374419744837, 0.07891374419744837, 138.436525793148, 1245.928732138332, 0.0290272375869813, 0.02902723758698131),
    @(3, 19.097892, 57.293676, 0.3678350062107417, 0.3678350062107418, 3, 37.82684066666667, 113.480522, 0.4118010609547572, 0.4118010609547572, 722.4129177532079, 6501.716259778872, 0.1514748458138832, 0.1514748458138832),
    @(3, 19.097892, 57.293676, 0.3678350062107417, 0.3678350062107418, 3, 38.20927633333334, 114.627829, 0.4159644383477588, 0.4159644383477588, 729.7166328121559, 6567.449695309404, 0.1530062817630956, 0.1530062817630956),
    @(3, 19.097892, 57.293676, 0.3678350062107417, 0.3678350062107418, 3, 8.572171666666666, 25.716515, 0.09332075650003555, 0.09332075650003555, 163.71040869546, 1473.39367825914, 0.03432664104678169, 0.03432664104678169),
    @(3, 0.144181, 0.432543, 0.002776998583428524, 0.002776998583428525, 3, 7.248785666666667, 21.746357, 0.07891374419744837, 0.07891374419744837, 1.045137166205667, 9.406234495851001, 0.0002191433558493551, 0.0002191433558493551),
    @(3, 0.144181, 0.432543, 0.002776998583428524, 0.002776998583428525, 3, 37.82684066666667, 113.480522, 0.4118010609547572, 0.4118010609547572, 5.453911714160666, 49.085205427446, 0.001143570962925724, 0.001143570962925724),
    @(3, 0.144181, 0.432543, 0.002776998583428524, 0.002776998583428525, 3, 38.20927633333334, 114.627829, 0.4159644383477588, 0.4159644383477588, 5.509051671016333, 49.581465039147, 0.001155132656048368, 0.001155132656048368),
    @(3, 0.144181, 0.432543, 0.002776998583428524, 0.002776998583428525, 3, 8.572171666666666, 25.716515, 0.09332075650003555, 0.09332075650003555, 1.235944283071666, 11.123498547645, 0.000259151608605077, 0.000259151608605077),
    @(3, 29.81579633333334, 89.44738900000002, 0.5742672348052801, 0.5742672348052801, 3, 7.248785666666667, 21.746357, 0.07891374419744837, 0.07891374419744837, 216.1283171013193, 1945.154853911874, 0.04531757766839989, 0.04531757766839989),
    @(3, 29.81579633333334, 89.44738900000002, 0.5742672348052801, 0.5742672348052801, 3, 37.82684066666667, 113.480522, 0.4118010609547572, 0.4118010609547572, 1127.837377250785, 10150.53639525706, 0.236483856564369, 0.236483856564369),
    @(3, 29.81579633333334, 89.44738900000002, 0.5742672348052801, 0.5742672348052801, 3, 38.20927633333334, 114.627829, 0.4159644383477588, 0.4159644383477588, 1139.24000119872, 10253.16001078848, 0.2388747477872989, 0.2388747477872989),
    @(3, 29.81579633333334, 89.44738900000002, 0.5742672348052801, 0.5742672348052801, 3, 8.572171666666666, 25.716515, 0.09332075650003555, 0.09332075650003555, 255.5861245477039, 2300.275120929335, 0.05359105278521228, 0.05359105278521228)
)

$startRow = 2
for ($i = 0; $i -lt $rowData.Count; $i++) {
    $r = $startRow + $i
    $vals = $rowData[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Range($cols[$j] + $r).Value = $vals[$j]
    }
}
